$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update Latest Handoff Datetime for the f15961c5 row (row 5)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-26 04:02:16"

# "de-de" sheet: update Latest Handoff Datetime for the f15961c5 row (row 5)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-26 04:02:25"
